# Member_import_template.xlsx edit:
# - Add an "Instructions" row at the top (row 1), merged A1:I1, with
#   explanatory text, bordered + wrapped styling, and a taller row height.
# - The old header row (row 1) becomes row 2; old sample-data row (row 2)
#   becomes row 3. Header row text turns bold gray; data row gets a plain
#   font.
# - Rename the "Street" header to "VotersID" and change the sample value
#   "JP Rizal" to "123-4455".
# - Freeze panes below the new header rows (top 2 rows frozen), with the
#   active selection on B4.
# - Drop the trailing 3 blank rows (940-942) that are no longer needed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Fix the two text values before shifting rows around (still on the
#    original row numbers here: header row = 1, sample row = 2).
# ---------------------------------------------------------------------
$ws.Range("H1").Value = "VotersID"
$ws.Range("H2").Value = "123-4455"

# ---------------------------------------------------------------------
# 2. Insert a new blank row at the top; existing rows 1.. shift to 2..
#    This also shifts the old blank placeholder row (was row 3) down to
#    row 4; since the sample row now permanently occupies row 3, that
#    extra blank placeholder is redundant, so remove it right away to
#    keep every later row number identical to the original file.
# ---------------------------------------------------------------------
$ws.Rows.Item(1).Insert()
$ws.Rows.Item(4).Delete()

# ---------------------------------------------------------------------
# 3. Populate the new instructions row.
# ---------------------------------------------------------------------
$instructions = "Instructions:`n1. Please read these instructions first to minimise potential errors during the data upload.`n2. Do not delete the row headings. Start from the row with the sample data (row 3).`n3. Email should be unique. If the email had already been saved before, the line will be ignored.`n4. Please fill in all data as correctly as possible`n5. Birthday should be in DATE format`n6. You can omit the RegionCode, ProvinceCode, and. CityCode. These will default to the Organisation" + [char]8217 + "s values, but can be changed later in the system.`n7. Barangay refers to the Barangay where the member" + [char]8217 + "s voting precinct is registered`n8. PositionInOrganisation refers to the Organisation you need to select in the dropdown when uploading the document in the system.`n9. IsRegisteredVoter should be either Y or N"

$ws.Range("A1").Value = $instructions
$ws.Range("A1:I1").Merge()
$ws.Rows.Item(1).RowHeight = 139.5

# Styling for the instructions cell: blue text, border, wrap, top-aligned.
$instrCell = $ws.Range("A1")
$instrCell.Font.Name = "Calibri"
$instrCell.Font.Bold = $false
$instrCell.Font.Color = 16711680
$instrCell.WrapText = $true
$instrCell.VerticalAlignment = -4160
$instrCell.Borders.Item(9).LineStyle = 1
$instrCell.Borders.Item(9).Weight = 2
$instrCell.Borders.Item(9).Color = 0

# The rest of row 1 (B1:Q1) just carries the border down, no content.
$restRow1 = $ws.Range("B1:Q1")
$restRow1.Borders.Item(9).LineStyle = 1
$restRow1.Borders.Item(9).Weight = 2
$restRow1.Borders.Item(9).Color = 0

# ---------------------------------------------------------------------
# 4. Re-style the header row (now row 2): bold, gray text.
# ---------------------------------------------------------------------
$headerRow = $ws.Range("A2:Q2")
$headerRow.Font.Color = 10066329

# ---------------------------------------------------------------------
# 5. Freeze panes under row 2, matching selection at B4.
# ---------------------------------------------------------------------
$ws.Range("A3").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("B4").Select()

# ---------------------------------------------------------------------
# 6. Remove the trailing blank rows 940-942 (the earlier insert/delete
#    pair above nets out to zero shift for every row from 4 onward, so
#    these are still the original row numbers).
# ---------------------------------------------------------------------
$ws.Range("A940:A942").EntireRow.Delete()
